$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-07-25 14:14:09"

for ($r = 261; $r -le 518; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
